$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New activity-log entry on row 3, mirroring the existing row 2 entry
# (same date / UBIT / duration), with a new description.

# Copy row 2's date cell formatting (keeps the existing built-in date style)
# then overwrite with the date value.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3").Value = 43864

$ws.Range("B3").Value = "jaclemon"
$ws.Range("C3").Value = "60 minutes"
$ws.Range("D3").Value = "Created function to allow program to intake command line arguments using C on notepad++"

$excel.CutCopyMode = 0

$ws.Range("D3").Select()
